$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19: H19: 1543.2222->1533.1111, I19: 1738.4->1720.2, K19: 1738.4->1720.2, M19: -1563.4->-1545.2
$ws.Range("H19").Value = 1533.1111
$ws.Range("I19").Value = 1720.2
$ws.Range("K19").Value = 1720.2
$ws.Range("M19").Value = -1545.2
# Row 62: H62: 7938->7333, J62: 7938->7333, L62: 7938->7333, N62: -9186->-8581
$ws.Range("H62").Value = 7333
$ws.Range("J62").Value = 7333
$ws.Range("L62").Value = 7333
$ws.Range("N62").Value = -8581
# Row 65: H65: 7938->7333, J65: 7938->7333, L65: 39690->36665, N65: -45930->-42905
$ws.Range("H65").Value = 7333
$ws.Range("J65").Value = 7333
$ws.Range("L65").Value = 36665
$ws.Range("N65").Value = -42905
# Row 88: H88: 4559->4465.7334, J88: 8027.875->7853, L88: 8027.875->7853, N88: -8839.875->-8665
$ws.Range("H88").Value = 4465.7334
$ws.Range("J88").Value = 7853
$ws.Range("L88").Value = 7853
$ws.Range("N88").Value = -8665
# Row 91: H91: 4559->4465.7334, J91: 8027.875->7853, L91: 8027.875->7853, N91: -10835.875->-10661
$ws.Range("H91").Value = 4465.7334
$ws.Range("J91").Value = 7853
$ws.Range("L91").Value = 7853
$ws.Range("N91").Value = -10661
# Row 99: H99: 1129.4445->1023, J99: 1755.6->1699.25, L99: 5266.799999999999->5097.75, N99: -8262.799999999999->-8093.75
$ws.Range("H99").Value = 1023
$ws.Range("J99").Value = 1699.25
$ws.Range("L99").Value = 5097.75
$ws.Range("N99").Value = -8093.75
# Row 138: H138: 2413.024->2367.15, J138: 2082.88->1974.3914, L138: 6248.64->5923.174199999999, N138: -16528.64->-16203.1742
$ws.Range("H138").Value = 2367.15
$ws.Range("J138").Value = 1974.3914
$ws.Range("L138").Value = 5923.174199999999
$ws.Range("N138").Value = -16203.1742

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32: 5505.706->5034.3335, I32: 5505.706->5034.3335, K32: 5505.706->5034.3335, M32: -5218.706->-4747.3335
$ws.Range("H32").Value = 5034.3335
$ws.Range("I32").Value = 5034.3335
$ws.Range("K32").Value = 5034.3335
$ws.Range("M32").Value = -4747.3335
# Row 45: H45: 4375.625->5538.75, I45: 5580->6533.3335, J45: 2368.3333->2555, K45: 5580->6533.3335, L45: 2368.3333->2555, M45: -5203->-6156.3335, N45: -3122.3333->-3309
$ws.Range("H45").Value = 5538.75
$ws.Range("I45").Value = 6533.3335
$ws.Range("J45").Value = 2555
$ws.Range("K45").Value = 6533.3335
$ws.Range("L45").Value = 2555
$ws.Range("M45").Value = -6156.3335
$ws.Range("N45").Value = -3309
# Row 61: H61: 3494.6->3397.9678, I61: 1735.9231->1683.1428, J61: 4839.4707->4810.1763, K61: 1735.9231->1683.1428, L61: 4839.4707->4810.1763, M61: -1523.9231->-1471.1428, N61: -5263.4707->-5234.1763
$ws.Range("H61").Value = 3397.9678
$ws.Range("I61").Value = 1683.1428
$ws.Range("J61").Value = 4810.1763
$ws.Range("K61").Value = 1683.1428
$ws.Range("L61").Value = 4810.1763
$ws.Range("M61").Value = -1471.1428
$ws.Range("N61").Value = -5234.1763
# Row 88: H88: 715.5->750, I88: 773->750, J88: 686.75->0, K88: 773->750, L88: 686.75->0, M88: -367->-344, N88: -1498.75->(removed)
$ws.Range("H88").Value = 750
$ws.Range("I88").Value = 750
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 750
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -344
$ws.Range("N88").ClearContents()
# Row 91: H91: 715.5->750, I91: 773->750, J91: 686.75->0, K91: 773->750, L91: 686.75->0, M91: 631->654, N91: -3494.75->(removed)
$ws.Range("H91").Value = 750
$ws.Range("I91").Value = 750
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 750
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 654
$ws.Range("N91").ClearContents()
# Row 97: H97: 997.5->915, I97: 1000->875, K97: 1000->875, M97: -504->-379
$ws.Range("H97").Value = 915
$ws.Range("I97").Value = 875
$ws.Range("K97").Value = 875
$ws.Range("M97").Value = -379
# Row 102: H102: 399->394, I102: 399->394, J102: 0->394, K102: 399->394, L102: 0->394, M102: 1223->1228, N102: (none)->-3638
$ws.Range("H102").Value = 394
$ws.Range("I102").Value = 394
$ws.Range("J102").Value = 394
$ws.Range("K102").Value = 394
$ws.Range("L102").Value = 394
$ws.Range("M102").Value = 1228
$ws.Range("N102").Value = -3638
# Row 122: H122: 1552->1388.75, I122: 1362.4->1235.1666, J122: 2500->1849.5, K122: 4087.2->3705.4998, L122: 7500->5548.5, M122: -1637.2->-1255.4998, N122: -12400->-10448.5
$ws.Range("H122").Value = 1388.75
$ws.Range("I122").Value = 1235.1666
$ws.Range("J122").Value = 1849.5
$ws.Range("K122").Value = 3705.4998
$ws.Range("L122").Value = 5548.5
$ws.Range("M122").Value = -1255.4998
$ws.Range("N122").Value = -10448.5
# Row 136: H136: 3494.6->3397.9678, I136: 1735.9231->1683.1428, J136: 4839.4707->4810.1763, K136: 5207.7693->5049.428400000001, L136: 14518.4121->14430.5289, M136: -2657.7693->-2499.428400000001, N136: -19618.4121->-19530.5289
$ws.Range("H136").Value = 3397.9678
$ws.Range("I136").Value = 1683.1428
$ws.Range("J136").Value = 4810.1763
$ws.Range("K136").Value = 5049.428400000001
$ws.Range("L136").Value = 14430.5289
$ws.Range("M136").Value = -2499.428400000001
$ws.Range("N136").Value = -19530.5289

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80: H80: 363.9524->352.2, J80: 197.61539->164.16667, L80: 197.61539->164.16667, N80: -2193.61539->-2160.16667
$ws.Range("H80").Value = 352.2
$ws.Range("J80").Value = 164.16667
$ws.Range("L80").Value = 164.16667
$ws.Range("N80").Value = -2160.16667
# Row 83: H83: 363.9524->352.2, J83: 197.61539->164.16667, L83: 988.0769499999999->820.8333500000001, N83: -10972.07695->-10804.83335
$ws.Range("H83").Value = 352.2
$ws.Range("J83").Value = 164.16667
$ws.Range("L83").Value = 820.8333500000001
$ws.Range("N83").Value = -10804.83335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31: 1893.0952->1928, I31: 1581.7693->1614, K31: 1581.7693->1614, M31: -1286.7693->-1319
$ws.Range("H31").Value = 1928
$ws.Range("I31").Value = 1614
$ws.Range("K31").Value = 1614
$ws.Range("M31").Value = -1319
# Row 34: H34: 1893.0952->1928, I34: 1581.7693->1614, K34: 1581.7693->1614, M34: -1379.7693->-1412
$ws.Range("H34").Value = 1928
$ws.Range("I34").Value = 1614
$ws.Range("K34").Value = 1614
$ws.Range("M34").Value = -1412
# Row 58: H58: 1475.375->1478.9231, I58: 1475.375->1478.9231, K58: 1475.375->1478.9231, M58: -1272.375->-1275.9231
$ws.Range("H58").Value = 1478.9231
$ws.Range("I58").Value = 1478.9231
$ws.Range("K58").Value = 1478.9231
$ws.Range("M58").Value = -1275.9231
# Row 132: H132: 2154.6667->1977.3334, I132: 1974.5454->1880, J132: 2650->2366.6667, K132: 5923.6362->5640, L132: 7950->7100.000100000001, M132: -3393.6362->-3110, N132: -13010->-12160.0001
$ws.Range("H132").Value = 1977.3334
$ws.Range("I132").Value = 1880
$ws.Range("J132").Value = 2366.6667
$ws.Range("K132").Value = 5640
$ws.Range("L132").Value = 7100.000100000001
$ws.Range("M132").Value = -3110
$ws.Range("N132").Value = -12160.0001
# Row 136: H136: 1475.375->1478.9231, I136: 1475.375->1478.9231, K136: 4426.125->4436.7693, M136: -1876.125->-1886.7693
$ws.Range("H136").Value = 1478.9231
$ws.Range("I136").Value = 1478.9231
$ws.Range("K136").Value = 4436.7693
$ws.Range("M136").Value = -1886.7693

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2: H2: 184404.25->158060.36, I2: 220127.8->220090.6, J2: 158887.42->123599.11, K2: 1320766.8->1320543.6, L2: 953324.52->741594.66, M2: -1320653.8->-1320430.6, N2: -953550.52->-741820.66
$ws.Range("H2").Value = 158060.36
$ws.Range("I2").Value = 220090.6
$ws.Range("J2").Value = 123599.11
$ws.Range("K2").Value = 1320543.6
$ws.Range("L2").Value = 741594.66
$ws.Range("M2").Value = -1320430.6
$ws.Range("N2").Value = -741820.66
# Row 46: H46: 2804->818, J46: 2804->818, L46: 8412->2454, N46: -8594->-2636
$ws.Range("H46").Value = 818
$ws.Range("J46").Value = 818
$ws.Range("L46").Value = 2454
$ws.Range("N46").Value = -2636
# Row 68: H68: 3914618.5->3914597.8, J68: 4473785->4473761.5, L68: 13421355->13421284.5, N68: -13422977->-13422906.5
$ws.Range("H68").Value = 3914597.8
$ws.Range("J68").Value = 4473761.5
$ws.Range("L68").Value = 13421284.5
$ws.Range("N68").Value = -13422906.5
# Row 71: H71: 3914618.5->3914597.8, J71: 4473785->4473761.5, L71: 40264065->40263853.5, N71: -40272177->-40271965.5
$ws.Range("H71").Value = 3914597.8
$ws.Range("J71").Value = 4473761.5
$ws.Range("L71").Value = 40263853.5
$ws.Range("N71").Value = -40271965.5
# Row 113: H113: 1363.5->2404.6, J113: 1221.1111->2344.5, L113: 3663.3333->7033.5, N113: -8003.3333->-11373.5
$ws.Range("H113").Value = 2404.6
$ws.Range("J113").Value = 2344.5
$ws.Range("L113").Value = 7033.5
$ws.Range("N113").Value = -11373.5
# Row 127: H127: 1987.25->1984, J127: 1987.25->1984, L127: 5961.75->5952, N127: -15881.75->-15872
$ws.Range("H127").Value = 1984
$ws.Range("J127").Value = 1984
$ws.Range("L127").Value = 5952
$ws.Range("N127").Value = -15872
# Row 131: H131: 2926.6191->2905.1765, I131: 0->1995, J131: 2926.6191->2962.0625, K131: 0->5985, L131: 8779.8573->8886.1875, M131: (none)->-945, N131: -18859.8573->-18966.1875
$ws.Range("H131").Value = 2905.1765
$ws.Range("I131").Value = 1995
$ws.Range("J131").Value = 2962.0625
$ws.Range("K131").Value = 5985
$ws.Range("L131").Value = 8886.1875
$ws.Range("M131").Value = -945
$ws.Range("N131").Value = -18966.1875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70: 0->9642.857, I70: 0->8750, J70: 0->10000, K70: 0->8750, L70: 0->10000, M70: (none)->-8480, N70: (none)->-10540
$ws.Range("H70").Value = 9642.857
$ws.Range("I70").Value = 8750
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 8750
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -8480
$ws.Range("N70").Value = -10540
# Row 73: H73: 0->9642.857, I73: 0->8750, J73: 0->10000, K73: 0->8750, L73: 0->10000, M73: (none)->-7814, N73: (none)->-11872
$ws.Range("H73").Value = 9642.857
$ws.Range("I73").Value = 8750
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 8750
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -7814
$ws.Range("N73").Value = -11872
# Row 80: H80: 1699.5->1633, J80: 0->1500, L80: 0->1500, N80: (none)->-3496
$ws.Range("H80").Value = 1633
$ws.Range("J80").Value = 1500
$ws.Range("L80").Value = 1500
$ws.Range("N80").Value = -3496
# Row 83: H83: 1699.5->1633, J83: 0->1500, L83: 0->7500, N83: (none)->-17484
$ws.Range("H83").Value = 1633
$ws.Range("J83").Value = 1500
$ws.Range("L83").Value = 7500
$ws.Range("N83").Value = -17484
# Row 113: H113: 4000->5000, I113: 4000->5000, K113: 4000->5000, M113: -1830->-2830
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
# Row 122: H122: 2194.6667->928, I122: 2194.6667->997.75, J122: 0->649, K122: 6584.000100000001->2993.25, L122: 0->1947, M122: -4134.000100000001->-543.25, N122: (none)->-6847
$ws.Range("H122").Value = 928
$ws.Range("I122").Value = 997.75
$ws.Range("J122").Value = 649
$ws.Range("K122").Value = 2993.25
$ws.Range("L122").Value = 1947
$ws.Range("M122").Value = -543.25
$ws.Range("N122").Value = -6847
# Row 132: H132: 1497.8->1497.6, I132: 1497->1496.6666, K132: 4491->4489.9998, M132: -1961->-1959.9998
$ws.Range("H132").Value = 1497.6
$ws.Range("I132").Value = 1496.6666
$ws.Range("K132").Value = 4489.9998
$ws.Range("M132").Value = -1959.9998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46: H46: 65170.125->2984.5557, I46: 74051.57000000001->3060.1667, J46: 3000->2833.3333, K46: 74051.57000000001->3060.1667, L46: 3000->2833.3333, M46: -73863.57000000001->-2872.1667, N46: -3376->-3209.3333
$ws.Range("H46").Value = 2984.5557
$ws.Range("I46").Value = 3060.1667
$ws.Range("J46").Value = 2833.3333
$ws.Range("K46").Value = 3060.1667
$ws.Range("L46").Value = 2833.3333
$ws.Range("M46").Value = -2872.1667
$ws.Range("N46").Value = -3209.3333
# Row 61: H61: 4465.8->3915.5, J61: 3443.3333->3220.8333, L61: 3443.3333->3220.8333, N61: -3847.3333->-3624.8333
$ws.Range("H61").Value = 3915.5
$ws.Range("J61").Value = 3220.8333
$ws.Range("L61").Value = 3220.8333
$ws.Range("N61").Value = -3624.8333
# Row 82: H82: 1109.4286->1095.4286, I82: 1109.4286->1111.1666, J82: 0->1001, K82: 1109.4286->1111.1666, L82: 0->1001, M82: -748.4286->-750.1666, N82: (none)->-1723
$ws.Range("H82").Value = 1095.4286
$ws.Range("I82").Value = 1111.1666
$ws.Range("J82").Value = 1001
$ws.Range("K82").Value = 1111.1666
$ws.Range("L82").Value = 1001
$ws.Range("M82").Value = -750.1666
$ws.Range("N82").Value = -1723
# Row 85: H85: 1109.4286->1095.4286, I85: 1109.4286->1111.1666, J85: 0->1001, K85: 1109.4286->1111.1666, L85: 0->1001, M85: 138.5714->136.8334, N85: (none)->-3497
$ws.Range("H85").Value = 1095.4286
$ws.Range("I85").Value = 1111.1666
$ws.Range("J85").Value = 1001
$ws.Range("K85").Value = 1111.1666
$ws.Range("L85").Value = 1001
$ws.Range("M85").Value = 136.8334
$ws.Range("N85").Value = -3497
# Row 93: H93: 1568.75->1991.6666, I93: 1428.3334->1992.5, K93: 1428.3334->1992.5, M93: -180.3334->-744.5
$ws.Range("H93").Value = 1991.6666
$ws.Range("I93").Value = 1992.5
$ws.Range("K93").Value = 1992.5
$ws.Range("M93").Value = -744.5
# Row 113: H113: 4465.8->3915.5, J113: 3443.3333->3220.8333, L113: 3443.3333->3220.8333, N113: -7783.3333->-7560.8333
$ws.Range("H113").Value = 3915.5
$ws.Range("J113").Value = 3220.8333
$ws.Range("L113").Value = 3220.8333
$ws.Range("N113").Value = -7560.8333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100: H100: 6252256->5557633, I100: 7144506.5->6668239.5, J100: 6500->4599.6665, K100: 14289013->13336479, L100: 13000->9199.333000000001, M100: -14288472->-13335938, N100: -14082->-10281.333
$ws.Range("H100").Value = 5557633
$ws.Range("I100").Value = 6668239.5
$ws.Range("J100").Value = 4599.6665
$ws.Range("K100").Value = 13336479
$ws.Range("L100").Value = 9199.333000000001
$ws.Range("M100").Value = -13335938
$ws.Range("N100").Value = -10281.333
# Row 122: H122: 1054.909->1158.7778, I122: 901->990.5714, K122: 2703->2971.7142, M122: -253->-521.7142000000003
$ws.Range("H122").Value = 1158.7778
$ws.Range("I122").Value = 990.5714
$ws.Range("K122").Value = 2971.7142
$ws.Range("M122").Value = -521.7142000000003
# Row 136: H136: 5099.6665->5074.25, J136: 4999->4998.5, L136: 14997->14995.5, N136: -20097->-20095.5
$ws.Range("H136").Value = 5074.25
$ws.Range("J136").Value = 4998.5
$ws.Range("L136").Value = 14995.5
$ws.Range("N136").Value = -20095.5
